# mod @20171226 by yxq
# Update the "支线高速耗时" (relay duration) and "高速网中转耗时"
# (transfer duration) parameter values on the Settings/Parameters sheet.
#   175分钟 -> 120分钟   (duration_relay rows)
#   15分钟  -> 10分钟    (duration_transfer rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Rows where column C = "duration_relay": E2, E8, E17 hold "175分钟"
$ws.Range("E2").Value  = "120分钟"
$ws.Range("E8").Value  = "120分钟"
$ws.Range("E17").Value = "120分钟"

# Rows where column C = "duration_transfer": E3, E9, E18 hold "15分钟"
$ws.Range("E3").Value  = "10分钟"
$ws.Range("E9").Value  = "10分钟"
$ws.Range("E18").Value = "10分钟"

# Reflect the author's last on-screen selection/scroll position.
$ws.Range("A2").Select()
$ws.Range("E17").Select()
